$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns to the updated, "reference column mapping" friendly labels.
# (Set in B,C,D,E,A order so new shared-string entries are appended in that order.)
$ws.Range("B1").Value = "Control ID"
$ws.Range("C1").Value = "Control Description"
$ws.Range("D1").Value = "Control Frequency"
$ws.Range("E1").Value = "Control Type"
$ws.Range("A1").Value = "Audit Leader From AE"

# Column A now holds a longer header, so re-fit its width to the content.
$ws.Columns.Item(1).AutoFit() | Out-Null

# Leave the selection where the user last clicked while reviewing the change.
$ws.Range("F4").Select() | Out-Null
